# Apply updated vm_pu results for "case with 380 kV done"
# Column B: slack/reference bus voltage set-point changed from 1.05 to 1.02 p.u.
# Columns C-F and I-N: recalculated bus voltage magnitudes (p.u.) for the new setup.
# Column G remains 1 (unchanged); column H has no data (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object 'object[,]' 24,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.077031357759747
$arrBF[0,2] = 1.063533384468024
$arrBF[0,3] = 1.077878684213367
$arrBF[0,4] = 1.082105847839071
$arrBF[1,0] = 1.02
$arrBF[1,1] = 1.07934777327441
$arrBF[1,2] = 1.064584468655403
$arrBF[1,3] = 1.079792744691825
$arrBF[1,4] = 1.083733136038729
$arrBF[2,0] = 1.02
$arrBF[2,1] = 1.080841386819635
$arrBF[2,2] = 1.065261125442498
$arrBF[2,3] = 1.081026506722317
$arrBF[2,4] = 1.084781520411548
$arrBF[3,0] = 1.02
$arrBF[3,1] = 1.081468075055466
$arrBF[3,2] = 1.065544773366595
$arrBF[3,3] = 1.081544065634548
$arrBF[3,4] = 1.085221184939216
$arrBF[4,0] = 1.02
$arrBF[4,1] = 1.081573227772099
$arrBF[4,2] = 1.065592351414172
$arrBF[4,3] = 1.081630901431111
$arrBF[4,4] = 1.085294944018316
$arrBF[5,0] = 1.02
$arrBF[5,1] = 1.080849765435179
$arrBF[5,2] = 1.06526491876156
$arrBF[5,3] = 1.081033426715069
$arrBF[5,4] = 1.084787399432194
$arrBF[6,0] = 1.02
$arrBF[6,1] = 1.077815310715613
$arrBF[6,2] = 1.063889327365322
$arrBF[6,3] = 1.078526551680056
$arrBF[6,4] = 1.082656758032253
$arrBF[7,0] = 1.02
$arrBF[7,1] = 1.072426423013238
$arrBF[7,2] = 1.061438325185423
$arrBF[7,3] = 1.074071476382956
$arrBF[7,4] = 1.078866304311075
$arrBF[8,0] = 1.02
$arrBF[8,1] = 1.068803635328532
$arrBF[8,2] = 1.059785443071059
$arrBF[8,3] = 1.071074480345075
$arrBF[8,4] = 1.076313844852727
$arrBF[9,0] = 1.02
$arrBF[9,1] = 1.067227294560549
$arrBF[9,2] = 1.059065089127202
$arrBF[9,3] = 1.06976998829416
$arrBF[9,4] = 1.07520226462541
$arrBF[10,0] = 1.02
$arrBF[10,1] = 1.066640581198432
$arrBF[10,2] = 1.058796805770543
$arrBF[10,3] = 1.069284392296428
$arrBF[10,4] = 1.074788395114358
$arrBF[11,0] = 1.02
$arrBF[11,1] = 1.066766487730114
$arrBF[11,2] = 1.058854385899807
$arrBF[11,3] = 1.069388602300915
$arrBF[11,4] = 1.074877216248143
$arrBF[12,0] = 1.02
$arrBF[12,1] = 1.067178821170237
$arrBF[12,2] = 1.059042927344022
$arrBF[12,3] = 1.069729870375452
$arrBF[12,4] = 1.07516807416449
$arrBF[13,0] = 1.02
$arrBF[13,1] = 1.067432714460423
$arrBF[13,2] = 1.059158999204235
$arrBF[13,3] = 1.069939996854857
$arrBF[13,4] = 1.075347150817485
$arrBF[14,0] = 1.02
$arrBF[14,1] = 1.06890808722432
$arrBF[14,2] = 1.059833151537998
$arrBF[14,3] = 1.071160909767801
$arrBF[14,4] = 1.076387480864249
$arrBF[15,0] = 1.02
$arrBF[15,1] = 1.069831472841084
$arrBF[15,2] = 1.060254776094671
$arrBF[15,3] = 1.071924920436448
$arrBF[15,4] = 1.077038335008255
$arrBF[16,0] = 1.02
$arrBF[16,1] = 1.070369332752339
$arrBF[16,2] = 1.060500255338869
$arrBF[16,3] = 1.072369903626467
$arrBF[16,4] = 1.077417356420475
$arrBF[17,0] = 1.02
$arrBF[17,1] = 1.070552605591695
$arrBF[17,2] = 1.060583881978257
$arrBF[17,3] = 1.072521521882372
$arrBF[17,4] = 1.077546490085442
$arrBF[18,0] = 1.02
$arrBF[18,1] = 1.069732478652605
$arrBF[18,2] = 1.060209586129048
$arrBF[18,3] = 1.071843016920087
$arrBF[18,4] = 1.076968567878198
$arrBF[19,0] = 1.02
$arrBF[19,1] = 1.067057432437017
$arrBF[19,2] = 1.058987426351298
$arrBF[19,3] = 1.069629404666253
$arrBF[19,4] = 1.075082450991791
$arrBF[20,0] = 1.02
$arrBF[20,1] = 1.065368617724247
$arrBF[20,2] = 1.058214880468314
$arrBF[20,3] = 1.068231530064613
$arrBF[20,4] = 1.073890896677326
$arrBF[21,0] = 1.02
$arrBF[21,1] = 1.066264558901111
$arrBF[21,2] = 1.058624817438113
$arrBF[21,3] = 1.068973157805413
$arrBF[21,4] = 1.074523108926362
$arrBF[22,0] = 1.02
$arrBF[22,1] = 1.069777212135691
$arrBF[22,2] = 1.060230006910929
$arrBF[22,3] = 1.071880027605005
$arrBF[22,4] = 1.077000094527715
$arrBF[23,0] = 1.02
$arrBF[23,1] = 1.07382474605383
$arrBF[23,2] = 1.06207524341828
$arrBF[23,3] = 1.075227847849791
$arrBF[23,4] = 1.079850624287458
$ws.Range("B2:F25").Value2 = $arrBF

$arrIN = New-Object 'object[,]' 24,6
$arrIN[0,0] = 1.048969482910763
$arrIN[0,1] = 1.081927849501584
$arrIN[0,2] = 1.066251766231477
$arrIN[0,3] = 1.080558826006627
$arrIN[0,4] = 1.084774930240696
$arrIN[0,5] = 1.083464311146735
$arrIN[1,0] = 1.049347159436503
$arrIN[1,1] = 1.083897997327971
$arrIN[1,2] = 1.067117287426222
$arrIN[1,3] = 1.082287992183687
$arrIN[1,4] = 1.086218833893026
$arrIN[1,5] = 1.085437256809015
$arrIN[2,0] = 1.049588200093771
$arrIN[2,1] = 1.085167444801169
$arrIN[2,2] = 1.067673276389993
$arrIN[2,3] = 1.083401694640658
$arrIN[2,4] = 1.087148074224778
$arrIN[2,5] = 1.086708507043232
$arrIN[3,0] = 1.049688739184866
$arrIN[3,1] = 1.085699862811964
$arrIN[3,2] = 1.067906052501598
$arrIN[3,3] = 1.083868678354168
$arrIN[3,4] = 1.087537534149164
$arrIN[3,5] = 1.087241681148672
$arrIN[4,0] = 1.049705573715659
$arrIN[4,1] = 1.085789185191269
$arrIN[4,2] = 1.067945080549178
$arrIN[4,3] = 1.083947016319646
$arrIN[4,4] = 1.087602856795116
$arrIN[4,5] = 1.087331130376001
$arrIN[5,0] = 1.049589546615965
$arrIN[5,1] = 1.085174563901113
$arrIN[5,2] = 1.067676390520022
$arrIN[5,3] = 1.083407939248495
$arrIN[5,4] = 1.087153282868065
$arrIN[5,5] = 1.086715636253114
$arrIN[6,0] = 1.049097815648422
$arrIN[6,1] = 1.082594798937264
$arrIN[6,2] = 1.066545120479757
$arrIN[6,3] = 1.081144293781077
$arrIN[6,4] = 1.08526396392115
$arrIN[6,5] = 1.084132207727114
$arrIN[7,0] = 1.048205481372919
$arrIN[7,1] = 1.07800653453221
$arrIN[7,2] = 1.0645200933573
$arrIN[7,3] = 1.077114664824047
$arrIN[7,4] = 1.081895107809791
$arrIN[7,5] = 1.079537427460324
$arrIN[8,0] = 1.04759286108149
$arrIN[8,1] = 1.074917412536218
$arrIN[8,2] = 1.063148157459192
$arrIN[8,3] = 1.074399271105366
$arrIN[8,4] = 1.079621331326787
$arrIN[8,5] = 1.07644391855677
$arrIN[9,0] = 1.047323303320898
$arrIN[9,1] = 1.073572206649346
$arrIN[9,2] = 1.062548741093807
$arrIN[9,3] = 1.073216255610803
$arrIN[9,4] = 1.078629874043384
$arrIN[9,5] = 1.075096802323242
$arrIN[10,0] = 1.04722252606108
$arrIN[10,1] = 1.073071361374436
$arrIN[10,2] = 1.062325272521142
$arrIN[10,3] = 1.072775714231519
$arrIN[10,4] = 1.078260542423144
$arrIN[10,5] = 1.074595245790593
$arrIN[11,0] = 1.047244172723021
$arrIN[11,1] = 1.073178848162634
$arrIN[11,2] = 1.062373244553463
$arrIN[11,3] = 1.072870262878524
$arrIN[11,4] = 1.078339813677796
$arrIN[11,5] = 1.07470288522236
$arrIN[12,0] = 1.047314986373321
$arrIN[12,1] = 1.073530830814107
$arrIN[12,2] = 1.062530285891046
$arrIN[12,3] = 1.073179863280964
$arrIN[12,4] = 1.078599366745593
$arrIN[12,5] = 1.07505536772957
$arrIN[13,0] = 1.047358530454471
$arrIN[13,1] = 1.073747542072964
$arrIN[13,2] = 1.062626935355266
$arrIN[13,3] = 1.073370469414072
$arrIN[13,4] = 1.078759144738539
$arrIN[13,5] = 1.075272386743272
$arrIN[14,0] = 1.047610659771307
$arrIN[14,1] = 1.075006526435744
$arrIN[14,2] = 1.063187824692715
$arrIN[14,3] = 1.074477629008624
$arrIN[14,4] = 1.079686983726247
$arrIN[14,5] = 1.076533159008255
$arrIN[15,0] = 1.047767660530562
$arrIN[15,1] = 1.075794196138429
$arrIN[15,2] = 1.063538211524667
$arrIN[15,3] = 1.075170163566775
$arrIN[15,4] = 1.080267129234235
$arrIN[15,5] = 1.077321947292265
$arrIN[16,0] = 1.047858823065878
$arrIN[16,1] = 1.076252900053
$arrIN[16,2] = 1.063742070088721
$arrIN[16,3] = 1.075573411531345
$arrIN[16,4] = 1.080604854471984
$arrIN[16,5] = 1.077781302619008
$arrIN[17,0] = 1.047889837236036
$arrIN[17,1] = 1.076409183371877
$arrIN[17,2] = 1.063811493448188
$arrIN[17,3] = 1.075710791565755
$arrIN[17,4] = 1.080719898175752
$arrIN[17,5] = 1.077937807878124
$arrIN[18,0] = 1.047750858640929
$arrIN[18,1] = 1.075709762442306
$arrIN[18,2] = 1.063500671818237
$arrIN[18,3] = 1.075095933342245
$arrIN[18,4] = 1.080204953946696
$arrIN[18,5] = 1.077237393690609
$arrIN[19,0] = 1.047294151560826
$arrIN[19,1] = 1.073427213439103
$arrIN[19,2] = 1.062484063839612
$arrIN[19,3] = 1.073088724751959
$arrIN[19,4] = 1.078522964299239
$arrIN[19,5] = 1.074951603206009
$arrIN[20,0] = 1.047003228659693
$arrIN[20,1] = 1.0719852631743
$arrIN[20,2] = 1.061840138749923
$arrIN[20,3] = 1.071820236531697
$arrIN[20,4] = 1.077459284213069
$arrIN[20,5] = 1.073507605206435
$arrIN[21,0] = 1.047157812444331
$arrIN[21,1] = 1.072750326609911
$arrIN[21,2] = 1.062181949808743
$arrIN[21,3] = 1.072493310256447
$arrIN[21,4] = 1.078023751897594
$arrIN[21,5] = 1.07427375511988
$arrIN[22,0] = 1.047758451968678
$arrIN[22,1] = 1.07574791665109
$arrIN[22,2] = 1.063517635990204
$arrIN[22,3] = 1.075129476934179
$arrIN[22,4] = 1.08023305033322
$arrIN[22,5] = 1.077275602082745
$arrIN[23,0] = 1.048439269789967
$arrIN[23,1] = 1.07919791559384
$arrIN[23,2] = 1.06504742381184
$arrIN[23,3] = 1.078161412647608
$arrIN[23,4] = 1.082770859564601
$arrIN[23,5] = 1.080730500419714
$ws.Range("I2:N25").Value2 = $arrIN

Write-Host "vm_pu.xlsx updated for case with 380 kV"
